$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain-looking number (e.g. "1.003") need an
# explicit Text number format first, otherwise Excel would store them as
# numeric values instead of preserving the original text representation
# used throughout this price list. Looping per-cell (instead of a single
# comma Union range) ensures the format is reliably applied to every cell.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D25", "D26", "D27", "D29", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.474.99"
$ws.Range("E2").Value = "  +3.41%  "
$ws.Range("D3").Value = "1.819.17"
$ws.Range("E3").Value = "  +4.94%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.77%  "
$ws.Range("D5").Value = "343.47"
$ws.Range("E5").Value = "  +2.85%  "
$ws.Range("D6").Value = "0.9998"
$ws.Range("D7").Value = "0.3840"
$ws.Range("E7").Value = "  +3.34%  "
$ws.Range("D8").Value = "0.3526"
$ws.Range("E8").Value = "  +4.12%  "
$ws.Range("D9").Value = "49.03"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").Value = "1.236"
$ws.Range("E10").Value = "  +2.52%  "
$ws.Range("D11").Value = "0.07810"
$ws.Range("E11").Value = "  +3.54%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").Value = "22.39"
$ws.Range("E13").Value = "  +9.13%  "
$ws.Range("D14").Value = "6.612"
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("D15").Value = "1.818.25"
$ws.Range("E15").Value = "  +5.78%  "
$ws.Range("D16").Value = "7.234"
$ws.Range("E16").Value = "  +3.20%  "
$ws.Range("D17").Value = "0.00001121"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "0.06726"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").Value = "86.50"
$ws.Range("E19").Value = "  +3.68%  "
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").Value = "17.67"
$ws.Range("E21").Value = "  +5.23%  "
$ws.Range("D22").Value = "6.569"
$ws.Range("E22").Value = "  +6.25%  "
$ws.Range("D23").Value = "13.22"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "27.485.54"
$ws.Range("E24").Value = "  +3.77%  "
$ws.Range("D25").Value = "2.461"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "2.690"
$ws.Range("E26").Value = "  +7.05%  "
$ws.Range("D27").Value = "22.50"
$ws.Range("E27").Value = "  +14.78%  "
$ws.Range("E28").Value = "  +2.29%  "
$ws.Range("D29").Value = "153.82"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("D30").Value = "2.020.45"
$ws.Range("E30").Value = "  +5.62%  "
$ws.Range("D31").Value = "136.68"
$ws.Range("E31").Value = "  +3.87%  "
$ws.Range("D32").Value = "6.368"
$ws.Range("E32").Value = "  +3.87%  "
$ws.Range("D33").Value = "4.066"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("E34").Value = "  +5.42%  "
$ws.Range("D35").Value = "0.08811"
$ws.Range("E35").Value = "  +2.43%  "
$ws.Range("D36").Value = "1.688"
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("D37").Value = "5.642"
$ws.Range("E37").Value = "  +3.12%  "
$ws.Range("D38").Value = "0.7028"
$ws.Range("E38").Value = "  +12.21%  "
$ws.Range("D39").Value = "0.2265"
$ws.Range("E39").Value = "  +4.13%  "
$ws.Range("D40").Value = "0.06493"
$ws.Range("E40").Value = "  +2.04%  "
$ws.Range("D41").Value = "0.02405"
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("D42").Value = "8.996"
$ws.Range("E42").Value = "  +3.13%  "
$ws.Range("D43").Value = "1.294"
$ws.Range("E43").Value = "  +4.39%  "
$ws.Range("D44").Value = "14.81"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("D45").Value = "0.6609"
$ws.Range("E45").Value = "  +9.02%  "
$ws.Range("D46").Value = "0.9996"
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").Value = "3.958"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").Value = "2.189"
$ws.Range("E48").Value = "  +5.78%  "
$ws.Range("D49").Value = "132.83"
$ws.Range("E49").Value = "  +2.48%  "
$ws.Range("D50").Value = "0.07339"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "80.68"
$ws.Range("E51").Value = "  +3.60%  "
